$wb = $excel.ActiveWorkbook
$wsKB = $wb.Worksheets.Item("KB")
$wsRef = $wb.Worksheets.Item("Reference")

# ---------------------------------------------------------------
# Reference sheet: re-lay the reference/notes block.
# Clear the old B3:C11 block first so nothing stale is left behind
# once the new layout (rows 2-18) is written.
# ---------------------------------------------------------------
$wsRef.Range("B3:C11").ClearContents()

$wsRef.Range("B2").Value = "For row 2 to 112"

$wsRef.Range("B3").Value = "https://www.acurite.com/blog/soil-moisture-guide-for-plants-and-vegetables.html"
$wsRef.Hyperlinks.Add($wsRef.Range("B3"), "https://www.acurite.com/blog/soil-moisture-guide-for-plants-and-vegetables.html")

$wsRef.Range("B13").Value = "Discretization"
$wsRef.Range("B13").Font.Bold = $true
$wsRef.Range("B13").Font.Underline = $true

$wsRef.Range("B6").Value = "https://www.almanac.com/plant-ph# "
$wsRef.Hyperlinks.Add($wsRef.Range("B6"), "https://www.almanac.com/plant-ph#", " ")

$wsRef.Range("B5").Value = "For row 113 to 133, however, it is pH data for the VEGATABLES, but it shows their moist levels instead as based on the previous KB website above told about the vegatables, all of them have a certain moist level."

$wsRef.Range("E7").Value = "Where it uses FOR ALL knowledge representation here."

$wsRef.Range("B14").Value = "very_low"
$wsRef.Range("C14").Value = "0 to 20 %"

$wsRef.Range("B15").Value = "low"
$wsRef.Range("C15").Value = "21 to 40 %"

$wsRef.Range("B16").Value = "mid"
$wsRef.Range("C16").Value = "41 to 60 %"

$wsRef.Range("B17").Value = "high"
$wsRef.Range("C17").Value = "61 to 80 %"

$wsRef.Range("B18").Value = "very_high"
$wsRef.Range("C18").Value = "81 to 100%"

# pageSetup
$wsRef.PageSetup.PaperSize = 9
$wsRef.PageSetup.Orientation = 1

# ---------------------------------------------------------------
# KB sheet: move the scroll/selection away from the old spot.
# ---------------------------------------------------------------
$wsKB.Activate()
$wsKB.Range("A142").Select()
$excel.ActiveWindow.ScrollRow = 109
$excel.ActiveWindow.ScrollColumn = 1

# ---------------------------------------------------------------
# Reference sheet becomes the active tab / selected sheet.
# ---------------------------------------------------------------
$wsRef.Activate()
$wsRef.Range("F14").Select()

Write-Host "done"
